$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "ECs" sending-cluster rows (old rows 2-4); remaining rows shift up
$ws.Range("A2:T4").EntireRow.Delete()

# Refresh the numeric columns (E:T) for the remaining rows with the recomputed TPM values
$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.06003666666666666"
$ws.Range("H2").Value = [double]"0.18011"
$ws.Range("I2").Value = [double]"0.007162610180657564"
$ws.Range("J2").Value = [double]"0.007162610180657565"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"28.25464766666667"
$ws.Range("N2").Value = [double]"84.763943"
$ws.Range("O2").Value = [double]"0.3168758800036845"
$ws.Range("P2").Value = [double]"0.3168758800036845"
$ws.Range("Q2").Value = [double]"1.696314863747778"
$ws.Range("R2").Value = [double]"15.26683377373"
$ws.Range("S2").Value = [double]"0.002269658404119215"
$ws.Range("T2").Value = [double]"0.002269658404119215"

$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.06003666666666666"
$ws.Range("H3").Value = [double]"0.18011"
$ws.Range("I3").Value = [double]"0.007162610180657564"
$ws.Range("J3").Value = [double]"0.007162610180657565"
$ws.Range("K3").Value = [double]"1"
$ws.Range("L3").Value = [double]"0.3333333333333333"
$ws.Range("M3").Value = [double]"0.04671833333333333"
$ws.Range("N3").Value = [double]"0.140155"
$ws.Range("O3").Value = [double]"0.0005239461189519747"
$ws.Range("P3").Value = [double]"0.0005239461189519747"
$ws.Range("Q3").Value = [double]"0.002804813005555555"
$ws.Range("R3").Value = [double]"0.02524331705"
$ws.Range("S3").Value = [double]"3.752821805721433e-06"
$ws.Range("T3").Value = [double]"3.752821805721434e-06"

$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.06003666666666666"
$ws.Range("H4").Value = [double]"0.18011"
$ws.Range("I4").Value = [double]"0.007162610180657564"
$ws.Range("J4").Value = [double]"0.007162610180657565"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"60.86492733333333"
$ws.Range("N4").Value = [double]"182.594782"
$ws.Range("O4").Value = [double]"0.6826001738773636"
$ws.Range("P4").Value = [double]"0.6826001738773636"
$ws.Range("Q4").Value = [double]"3.654127354002222"
$ws.Range("R4").Value = [double]"32.88714618602"
$ws.Range("S4").Value = [double]"0.004889198954732628"
$ws.Range("T4").Value = [double]"0.004889198954732629"

$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"8.321917"
$ws.Range("H5").Value = [double]"24.965751"
$ws.Range("I5").Value = [double]"0.9928373898193424"
$ws.Range("J5").Value = [double]"0.9928373898193424"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"28.25464766666667"
$ws.Range("N5").Value = [double]"84.763943"
$ws.Range("O5").Value = [double]"0.3168758800036845"
$ws.Range("P5").Value = [double]"0.3168758800036845"
$ws.Range("Q5").Value = [double]"235.1328327462437"
$ws.Range("R5").Value = [double]"2116.195494716193"
$ws.Range("S5").Value = [double]"0.3146062215995653"
$ws.Range("T5").Value = [double]"0.3146062215995653"

$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"8.321917"
$ws.Range("H6").Value = [double]"24.965751"
$ws.Range("I6").Value = [double]"0.9928373898193424"
$ws.Range("J6").Value = [double]"0.9928373898193424"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.3333333333333333"
$ws.Range("M6").Value = [double]"0.04671833333333333"
$ws.Range("N6").Value = [double]"0.140155"
$ws.Range("O6").Value = [double]"0.0005239461189519747"
$ws.Range("P6").Value = [double]"0.0005239461189519747"
$ws.Range("Q6").Value = [double]"0.3887860923783333"
$ws.Range("R6").Value = [double]"3.499074831405"
$ws.Range("S6").Value = [double]"0.0005201932971462533"
$ws.Range("T6").Value = [double]"0.0005201932971462533"

$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"8.321917"
$ws.Range("H7").Value = [double]"24.965751"
$ws.Range("I7").Value = [double]"0.9928373898193424"
$ws.Range("J7").Value = [double]"0.9928373898193424"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"60.86492733333333"
$ws.Range("N7").Value = [double]"182.594782"
$ws.Range("O7").Value = [double]"0.6826001738773636"
$ws.Range("P7").Value = [double]"0.6826001738773636"
$ws.Range("Q7").Value = [double]"506.5128734790313"
$ws.Range("R7").Value = [double]"4558.615861311282"
$ws.Range("S7").Value = [double]"0.677710974922631"
$ws.Range("T7").Value = [double]"0.677710974922631"
